# Inventory.xlsx edit: remove the ASUS ExpertBook / Shrouq Aldakkan record
# (row 58) and widen the columns for readability.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete inventory row entirely - everything below shifts up.
$ws.Rows(58).Delete()

# Improve the layout by giving every column an explicit, generous width.
# (ColumnWidth is in "characters"; Excel stores width = ColumnWidth + 0.8333
# in the underlying <col> element, so subtract that padding to land on the
# exact target widths.)
$pad = 0.8333333333333334
$widths = @(15, 18, 133, 57, 16, 10, 10, 21, 29, 19, 31, 46)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $widths[$i] - $pad
}
